# Update NATMI LR-pair TPM-derived metrics (Ifna11-Ifnar2) with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.7994676808118871
$ws.Range("J2").Value = 0.7994676808118871
$ws.Range("M2").Value = 13.11310966666667
$ws.Range("N2").Value = 39.339329
$ws.Range("O2").Value = 0.2602886552498481
$ws.Range("P2").Value = 0.2602886552498481
$ws.Range("Q2").Value = 2.611427708714555
$ws.Range("R2").Value = 23.502849378431
$ws.Range("S2").Value = 0.2080923675542409
$ws.Range("T2").Value = 0.2080923675542409
$ws.Range("I3").Value = 0.7994676808118871
$ws.Range("J3").Value = 0.7994676808118871
$ws.Range("O3").Value = 0.5256625072946081
$ws.Range("P3").Value = 0.5256625072946082
$ws.Range("S3").Value = 0.420250185596582
$ws.Range("T3").Value = 0.4202501855965821
$ws.Range("I4").Value = 0.7994676808118871
$ws.Range("J4").Value = 0.7994676808118871
$ws.Range("M4").Value = 10.783589
$ws.Range("N4").Value = 32.350767
$ws.Range("O4").Value = 0.2140488374555438
$ws.Range("P4").Value = 0.2140488374555438
$ws.Range("Q4").Value = 2.147512209523667
$ws.Range("R4").Value = 19.327609885713
$ws.Range("S4").Value = 0.1711251276610642
$ws.Range("T4").Value = 0.1711251276610642
$ws.Range("G5").Value = 0.04995233333333333
$ws.Range("H5").Value = 0.149857
$ws.Range("I5").Value = 0.2005323191881128
$ws.Range("J5").Value = 0.2005323191881128
$ws.Range("M5").Value = 13.11310966666667
$ws.Range("N5").Value = 39.339329
$ws.Range("O5").Value = 0.2602886552498481
$ws.Range("P5").Value = 0.2602886552498481
$ws.Range("Q5").Value = 0.6550304251058888
$ws.Range("R5").Value = 5.895273825953
$ws.Range("S5").Value = 0.0521962876956072
$ws.Range("T5").Value = 0.0521962876956072
$ws.Range("G6").Value = 0.04995233333333333
$ws.Range("H6").Value = 0.149857
$ws.Range("I6").Value = 0.2005323191881128
$ws.Range("J6").Value = 0.2005323191881128
$ws.Range("O6").Value = 0.5256625072946081
$ws.Range("P6").Value = 0.5256625072946082
$ws.Range("Q6").Value = 1.322858021933
$ws.Range("R6").Value = 11.905722197397
$ws.Range("S6").Value = 0.105412321698026
$ws.Range("T6").Value = 0.1054123216980261
$ws.Range("G7").Value = 0.04995233333333333
$ws.Range("H7").Value = 0.149857
$ws.Range("I7").Value = 0.2005323191881128
$ws.Range("J7").Value = 0.2005323191881128
$ws.Range("M7").Value = 10.783589
$ws.Range("N7").Value = 32.350767
$ws.Range("O7").Value = 0.2140488374555438
$ws.Range("P7").Value = 0.2140488374555438
$ws.Range("Q7").Value = 0.5386654322576666
$ws.Range("R7").Value = 4.847988890318999
$ws.Range("S7").Value = 0.04292370979447959
$ws.Range("T7").Value = 0.0429237097944796
Write-Output "Updated TPM-derived values for rows 2-7"
